# "as on 30th june"
# Insert a new "adminuser" worksheet between "managenews" and "managecontact",
# populate it with a small admin-user record, and make it the active sheet.

$wb = $excel.ActiveWorkbook

$afterSheet = $wb.Worksheets.Item("managenews")
$newSheet = $wb.Worksheets.Add($null, $afterSheet)
$newSheet.Name = "adminuser"

# Fill A3 first, then A1, then A2 so the shared-string table and cell values
# line up with the authored workbook.
$newSheet.Range("A3").Value = "test3"
$newSheet.Range("A1").Value = "Ezabella"
$newSheet.Range("A2").Value = "ezabella"

# Leave the selection on A3, as the saved workbook shows.
$newSheet.Range("A3").Select() | Out-Null
